$d = $word.ActiveDocument

$pairs = @(
    @("492×9=4428", "985×8=7880"),
    @("246×8=1968", "849×8=6792"),
    @("829×8=6632", "581×6=3486"),
    @("282×6=1692", "309×7=2163"),
    @("342×5=1710", "830×7=5810"),
    @("909×4=3636", "224×2=448"),
    @("648×8=5184", "574×3=1722"),
    @("285×2=570", "863×2=1726"),
    @("623×9=5607", "371×9=3339"),
    @("349×4=1396", "183×2=366"),
    @("233×2=466", "139×3=417"),
    @("622×2=1244", "554×3=1662"),
    @("660×6=3960", "985×2=1970"),
    @("374×5=1870", "520×3=1560"),
    @("364×6=2184", "778×7=5446"),
    @("529×8=4232", "824×3=2472"),
    @("420×3=1260", "981×3=2943"),
    @("932×9=8388", "951×2=1902"),
    @("694×9=6246", "612×7=4284"),
    @("736×7=5152", "200×7=1400"),
    @("465×3=1395", "622×7=4354"),
    @("274×2=548", "375×6=2250"),
    @("122×8=976", "436×4=1744"),
    @("197×6=1182", "188×9=1692"),
    @("529×9=4761", "631×5=3155")
)

foreach ($pair in $pairs) {
    $old = $pair[0]
    $new = $pair[1]
    $range = $d.Content
    $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

$d.Save()
